$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B/C columns (coin name/link) as plain text, D/E as explicit text format
# to preserve values like "67.053.79" or "  -0.65%  " exactly as text (not numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.053.79'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.605.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.21'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.79%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.606.69'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.121'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.75%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.64%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.53'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.081.48'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.964.26'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.605.80'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '364.62'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.02'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.35'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.35%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.06'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.34'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.38%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.99'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.733.73'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '581.69'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.37'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.66'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.80'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -8.44%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.10'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.60%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.22'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.25%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '154.14'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0293'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.70'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.54%  '
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.69'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.617'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.85%  '
